$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K2").Value = 1.925925925925943
$ws.Range("R2").Value = 1.599224389494701
$ws.Range("S2").Value = 1.721244347680456

$ws.Range("K3").Value = 1.925925925925943
$ws.Range("R3").Value = 1.459904774678112
$ws.Range("S3").Value = 1.554373915558126

$ws.Range("K8").Value = 1.925925925925943
$ws.Range("R8").Value = 1.599224389494701
$ws.Range("S8").Value = 1.721244347680456

$ws.Range("K9").Value = 1.925925925925943
$ws.Range("R9").Value = 1.459904774678112
$ws.Range("S9").Value = 1.554373915558126

$ws.Range("K16").Value = 13.17361111111111
$ws.Range("R16").Value = 1.763755319824684
$ws.Range("S16").Value = 1.916157449486122

$ws.Range("K17").Value = 13.17361111111111
$ws.Range("R17").Value = 1.58937742977605
$ws.Range("S17").Value = 1.704024252511443

$ws.Range("K18").Value = 13.17361111111111
$ws.Range("R18").Value = 1.763755319824684
$ws.Range("S18").Value = 1.916157449486122

$ws.Range("K19").Value = 13.17361111111111
$ws.Range("R19").Value = 1.58937742977605
$ws.Range("S19").Value = 1.704024252511443
